$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Only evaluate/populate the relevant fields in the sheet (prevents OOM)
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1

$ws.Range("G5").Value = 1
$ws.Range("I5").Value = 1

$ws.Activate()
$ws.Range("I17").Select()
